# Update the LR-pairs worksheet with refreshed TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Egf/Erbb4 -> MuSCs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.041452
$ws.Range("H2").Value = 0.124356
$ws.Range("I2").Value = 0.05439747478414846
$ws.Range("J2").Value = 0.05439747478414846
$ws.Range("Q2").Value = 0.0002632340173333333
$ws.Range("R2").Value = 0.002369106156
$ws.Range("S2").Value = 0.05439747478414846
$ws.Range("T2").Value = 0.05439747478414846

# Row 3 (FAPs -> Egf/Erbb4 -> MuSCs)
$ws.Range("I3").Value = 0.5204718857143857
$ws.Range("J3").Value = 0.5204718857143856
$ws.Range("S3").Value = 0.5204718857143857
$ws.Range("T3").Value = 0.5204718857143856

# Row 4 (MuSCs -> Egf/Erbb4 -> MuSCs)
$ws.Range("G4").Value = 0.3239583333333333
$ws.Range("H4").Value = 0.971875
$ws.Range("I4").Value = 0.4251306395014658
$ws.Range("J4").Value = 0.4251306395014658
$ws.Range("Q4").Value = 0.002057243402777778
$ws.Range("R4").Value = 0.018515190625
$ws.Range("S4").Value = 0.4251306395014658
$ws.Range("T4").Value = 0.4251306395014658
